$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 149, shifting all existing data
# (rows 149-182) down to rows 151-184.
$ws.Rows.Item(149).Insert()
$ws.Rows.Item(149).Insert()

# New row 149 data
$ws.Cells.Item(149, 1).Value = 10
$ws.Cells.Item(149, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(149, 3).Value = "La Araucanía"
$ws.Cells.Item(149, 4).Value = "2022-08-25"
$ws.Cells.Item(149, 5).Value = 9
$ws.Cells.Item(149, 6).Value = 100112012
$ws.Cells.Item(149, 7).Value = "Espinaca"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 80
$ws.Cells.Item(149, 11).Value = 9000
$ws.Cells.Item(149, 12).Value = 9000
$ws.Cells.Item(149, 13).Value = 9000
$ws.Cells.Item(149, 14).Value = "`$/docena de atados"
$ws.Cells.Item(149, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(149, 16).Value = 3000
$ws.Cells.Item(149, 17).Value = 3
$ws.Cells.Item(149, 18).Value = "Hortaliza"

# New row 150 data
$ws.Cells.Item(150, 1).Value = 10
$ws.Cells.Item(150, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value = "La Araucanía"
$ws.Cells.Item(150, 4).Value = "2022-08-25"
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = 100112012
$ws.Cells.Item(150, 7).Value = "Espinaca"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 40
$ws.Cells.Item(150, 11).Value = 12000
$ws.Cells.Item(150, 12).Value = 12000
$ws.Cells.Item(150, 13).Value = 12000
$ws.Cells.Item(150, 14).Value = "`$/docena de atados"
$ws.Cells.Item(150, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(150, 16).Value = 4000
$ws.Cells.Item(150, 17).Value = 3
$ws.Cells.Item(150, 18).Value = "Hortaliza"
